# Update the "K" column (G) in the save_data sheet with recalculated strikeout
# values (regenerated using K instead of Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0, 1, 1, 2, 2, 2, 3, 1, 1, 3, 2, 4, 2, 1, 4, 1, 2, 0, 5, 7, 6, 5, 8, 3, 4, 2, 4, 4, 4, 4, 2, 6, 8, 3, 7, 2, 4, 10, 3, 7, 3, 3, 4, 2, 2, 2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
